# Capitalize the first letter of each value in column AB ("Specimen") for
# data rows (row 2 through the last used row), leaving the header row (row 1)
# untouched. This mirrors the commit message: "update isolate column
# comment if the isolate is not a clinical isolate." (values such as "NA"
# are already capitalized and are unaffected.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, "AB").End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Range("AB$r")
    $val = $cell.Value()
    if ($val -ne $null -and $val -is [string] -and $val.Length -gt 0) {
        $newVal = $val.Substring(0,1).ToUpper() + $val.Substring(1)
        $cell.Value = $newVal
    }
}
